# Resident Evil 4 report-card template: drop the "${imagePic}" merge field.
#
# The original "Image:" paragraph in the table cell read:
#     Image: ${imagePic}
# and was preceded by an empty, underlined spacer paragraph.
#
# The edit removes the spacer paragraph entirely (its paragraph mark is
# deleted, merging it into the following paragraph, whose own mark/pPr
# survives) and clears all the run content of the "Image:" paragraph,
# leaving it as a single empty paragraph (its pPr - sz 24 / szCs 24 /
# lang en-US - is kept intact).

$d = $word.ActiveDocument

# Locate the paragraph that still carries the "${imagePic}" merge field.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "imagePic") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    # The paragraph immediately before it is the empty, underlined spacer
    # paragraph - delete it outright (Range.Delete on a paragraph's Range,
    # which includes its end-of-paragraph mark, merges it into the next
    # paragraph and keeps that next paragraph's own mark/formatting).
    $prev = $target.Previous()
    if ($prev -ne $null -and $prev.Range.Text.Trim() -eq "") {
        $prev.Range.Delete()
    }

    # Re-resolve the target paragraph since the delete above shifted ranges.
    $target2 = $null
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -match "imagePic") {
            $target2 = $p
            break
        }
    }

    # Clear just the run text, excluding the trailing paragraph mark, so
    # the paragraph itself (and its pPr) survives as an empty paragraph.
    $r = $target2.Range
    $r2 = $d.Range($r.Start, $r.End - 1)
    $r2.Text = ""
}
